# Append 4 new arrival rows (310-313) to the "Main Data" sheet, mirroring
# the existing row layout:
#   A=NUMBER  B=DATE  C=TIME  D=FLIGHT  E=FROM  F=SHORT  G=AIRLINE
#   H=MODEL   I=AIRCFAT ID  J=STATUS  K=(blank)  L=DIFFERENCE  M=(blank)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 310 -> FR5669 EDI (EI-DHP)
$ws.Cells.Item(310, 1).Value = 309
$ws.Cells.Item(310, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(310, 3).Value = "4:50 PM"
$ws.Cells.Item(310, 4).Value = "FR5669"
$ws.Cells.Item(310, 5).Value = "Edinburgh"
$ws.Cells.Item(310, 6).Value = "(EDI)"
$ws.Cells.Item(310, 7).Value = "Ryanair "
$ws.Cells.Item(310, 8).Value = "B738"
$ws.Cells.Item(310, 9).Value = "(EI-DHP)"
$ws.Cells.Item(310, 10).Value = "4:33 PM"
$ws.Cells.Item(310, 12).Value = "0 hours, -17 minutes"

# Row 311 -> FR7101 OSL (SP-RSH)
$ws.Cells.Item(311, 1).Value = 310
$ws.Cells.Item(311, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(311, 3).Value = "5:15 PM"
$ws.Cells.Item(311, 4).Value = "FR7101"
$ws.Cells.Item(311, 5).Value = "Oslo"
$ws.Cells.Item(311, 6).Value = "(OSL)"
$ws.Cells.Item(311, 7).Value = "Ryanair "
$ws.Cells.Item(311, 8).Value = "B738"
$ws.Cells.Item(311, 9).Value = "(SP-RSH)"
$ws.Cells.Item(311, 10).Value = "5:12 PM"
$ws.Cells.Item(311, 12).Value = "0 hours, -3 minutes"

# Row 312 -> W61072 EIN (HA-LTC)
$ws.Cells.Item(312, 1).Value = 311
$ws.Cells.Item(312, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(312, 3).Value = "5:15 PM"
$ws.Cells.Item(312, 4).Value = "W61072"
$ws.Cells.Item(312, 5).Value = "Eindhoven"
$ws.Cells.Item(312, 6).Value = "(EIN)"
$ws.Cells.Item(312, 7).Value = "Wizz Air "
$ws.Cells.Item(312, 8).Value = "A321"
$ws.Cells.Item(312, 9).Value = "(HA-LTC)"
$ws.Cells.Item(312, 10).Value = "4:52 PM"
$ws.Cells.Item(312, 12).Value = "0 hours, -23 minutes"

# Row 313 -> FR6388 ATH (SP-RSB)
$ws.Cells.Item(313, 1).Value = 312
$ws.Cells.Item(313, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(313, 3).Value = "5:20 PM"
$ws.Cells.Item(313, 4).Value = "FR6388"
$ws.Cells.Item(313, 5).Value = "Athens"
$ws.Cells.Item(313, 6).Value = "(ATH)"
$ws.Cells.Item(313, 7).Value = "Ryanair "
$ws.Cells.Item(313, 8).Value = "B738"
$ws.Cells.Item(313, 9).Value = "(SP-RSB)"
$ws.Cells.Item(313, 10).Value = "5:14 PM"
$ws.Cells.Item(313, 12).Value = "0 hours, -6 minutes"
